$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F10").Value = "ppe"
$ws.Range("F16").Value = "application instructions"
$ws.Range("F17").Value = "application instructions"
$ws.Range("F18").Value = "application instructions"
$ws.Range("F19").Value = "application instructions"
$ws.Range("F39").Value = "application instructions"
$ws.Range("F40").Value = "application instructions"
$ws.Range("F41").Value = "application instructions"
